$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update the ticket price total (D24) which drives the shared formula
# in D2:D21 ($D$24/$A$21) from 1.35 to 1.5.
$ws.Range("D24").Value = 30

# N5 used to be a helper formula (=N2+N3); put it back to a plain
# literal value of 3.
$ws.Range("N5").Value = 3

# Update E22 ticket price from 3 to 5 (E24 = SUM(E2:E22) recalculates
# automatically).
$ws.Range("E22").Value = 5

# Move the active selection to N20.
$ws.Activate()
$ws.Range("N20").Select()
